$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.06992366666666666
$ws.Range("H2").Value = 0.209771
$ws.Range("M2").Value = 0.8059226666666667
$ws.Range("N2").Value = 2.417768
$ws.Range("O2").Value = 0.1314814101815314
$ws.Range("P2").Value = 0.1314814101815314
$ws.Range("Q2").Value = 0.05635306790311111
$ws.Range("R2").Value = 0.507177611128
$ws.Range("S2").Value = 0.1314814101815314
$ws.Range("T2").Value = 0.1314814101815314

# Row 3 updates
$ws.Range("G3").Value = 0.06992366666666666
$ws.Range("H3").Value = 0.209771
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.323633333333333
$ws.Range("N3").Value = 15.9709
$ws.Range("O3").Value = 0.8685185898184686
$ws.Range("P3").Value = 0.8685185898184687
$ws.Range("Q3").Value = 0.3722479626555555
$ws.Range("R3").Value = 3.3502316639
$ws.Range("S3").Value = 0.8685185898184686
$ws.Range("T3").Value = 0.8685185898184687
